$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 05:42"

# --- Row 24: Pakistan - refreshed case counts ---
$ws.Range("B24").Value = 313984
$ws.Range("C24").Value = 553
$ws.Range("D24").Value = 298593
$ws.Range("E24").Value = 8884
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 6507

# --- Row 36: Belgica - refreshed case counts ---
$ws.Range("B36").Value = 124234
$ws.Range("C36").Value = 3175
$ws.Range("D36").Value = 19521
$ws.Range("E36").Value = 94676
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 10037

# --- Row 39: Kazajistan - refreshed case counts ---
$ws.Range("B39").Value = 108106
$ws.Range("C39").Value = 62
$ws.Range("D39").Value = 103110
$ws.Range("E39").Value = 3271
$ws.Range("F39").Value = 0

# --- Rows 51-53: Honduras overtakes Chequia and Costa Rica in case counts ---
# Row 51 becomes Honduras with refreshed counts
$ws.Range("A51").Value = "Honduras"
$ws.Range("B51").Value = 78269
$ws.Range("C51").Value = 671
$ws.Range("D51").Value = 28978
$ws.Range("E51").Value = 46905
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 6
$ws.Range("H51").Value = 2386

# Row 52 becomes Chequia (previous row-51 values)
$ws.Range("A52").Value = "Chequia"
$ws.Range("B52").Value = 78051
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 35032
$ws.Range("E52").Value = 42320
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 699

# Row 53 becomes Costa Rica (previous row-52 values)
$ws.Range("A53").Value = "Costa Rica"
$ws.Range("B53").Value = 77829
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 42621
$ws.Range("E53").Value = 34278
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 930

# --- Row 56: Venezuela - refreshed case counts ---
$ws.Range("B56").Value = 76820
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 67216
$ws.Range("E56").Value = 8961
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 643

# --- Row 172: Islas Turcas y Caicos - refreshed case counts ---
$ws.Range("B172").Value = 695
$ws.Range("C172").Value = 5
$ws.Range("D172").Value = 651
$ws.Range("E172").Value = 38
$ws.Range("F172").Value = 0

# --- Row 173: San Martin (Parte Holandesa) - refreshed case counts ---
$ws.Range("B173").Value = 674
$ws.Range("C173").Value = 6
$ws.Range("D173").Value = 572
$ws.Range("E173").Value = 80
$ws.Range("F173").Value = 0
